$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply odds/value corrections to row 5 (Spartak Varna vs CSKA 1948 Sofia)

# Row 5 updates
$ws.Range("G5").Value = 2.4
$ws.Range("I5").Value = 3.1
$ws.Range("J5").Value = 3.1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 3.75
$ws.Range("S5").Value = 1.5
$ws.Range("T5").Value = 2.5
$ws.Range("U5").Value = 1.91
$ws.Range("V5").Value = 1.8
$ws.Range("W5").Value = 7
$ws.Range("AC5").Value = 8
$ws.Range("AG5").Value = 351
$ws.Range("AJ5").Value = 12
$ws.Range("AK5").Value = 34
$ws.Range("AM5").Value = 41
$ws.Range("AO5").Value = 13
$ws.Range("AT5").Value = 2.5
$ws.Range("AV5").Value = 67
$ws.Range("AX5").Value = 5
$ws.Range("BC5").Value = 251

# Row 6 updates
$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 9
$ws.Range("R6").Value = 1.62

# Row 7 updates
$ws.Range("Q7").Value = 1.57

# Rows 9-12: match fixtures changed (new matches replacing prior entries, plus a newly
# added fixture in row 12) - update every affected cell to the new fixture/odds data.

# Row 9 updates
$ws.Range("A9").Value = "UPa7e1IO"
$ws.Range("C9").Value = "11:25"
$ws.Range("D9").Value = "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE"
$ws.Range("E9").Value = "Al Qadisiya"
$ws.Range("F9").Value = "Al Khaleej"
$ws.Range("G9").Value = 1.57
$ws.Range("H9").Value = 4
$ws.Range("I9").Value = 4.75
$ws.Range("J9").Value = 2.1
$ws.Range("K9").Value = 2.2
$ws.Range("L9").Value = 5.5
$ws.Range("M9").Value = 1.02
$ws.Range("N9").Value = 11
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 3.5
$ws.Range("Q9").Value = 1.9
$ws.Range("R9").Value = 1.9
$ws.Range("S9").Value = 1.4
$ws.Range("T9").Value = 2.75
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = 1.73
$ws.Range("X9").Value = 7.5
$ws.Range("Z9").Value = 11
$ws.Range("AA9").Value = 13
$ws.Range("AC9").Value = 11
$ws.Range("AD9").Value = 8
$ws.Range("AE9").Value = 19
$ws.Range("AG9").Value = 800
$ws.Range("AH9").Value = 13
$ws.Range("AI9").Value = 26
$ws.Range("AJ9").Value = 17
$ws.Range("AK9").Value = 51
$ws.Range("AL9").Value = 41
$ws.Range("AN9").Value = 3.6
$ws.Range("AO9").Value = 8
$ws.Range("AP9").Value = 21
$ws.Range("AQ9").Value = 26
$ws.Range("AR9").Value = 51
$ws.Range("AS9").Value = 151
$ws.Range("AT9").Value = 2.75
$ws.Range("AU9").Value = 9
$ws.Range("AV9").Value = 51
$ws.Range("AW9").Value = 81
$ws.Range("AX9").Value = 7
$ws.Range("AY9").Value = 29
$ws.Range("BA9").Value = 101
$ws.Range("BB9").Value = 126
$ws.Range("BC9").Value = 500
$ws.Range("BD9").Value = 81

# Row 10 updates
$ws.Range("A10").Value = "rTPtA0XH"
$ws.Range("C10").Value = "11:40"
$ws.Range("E10").Value = "Al Nassr"
$ws.Range("F10").Value = "Damac"
$ws.Range("G10").Value = 1.25
$ws.Range("H10").Value = 6.25
$ws.Range("I10").Value = 8
$ws.Range("J10").Value = 1.62
$ws.Range("K10").Value = 2.75
$ws.Range("L10").Value = 7.5
$ws.Range("N10").Value = 12
$ws.Range("O10").Value = 1.13
$ws.Range("P10").Value = 5.5
$ws.Range("Q10").Value = 1.4
$ws.Range("R10").Value = 2.75
$ws.Range("S10").Value = 1.22
$ws.Range("T10").Value = 4
$ws.Range("U10").Value = 1.83
$ws.Range("V10").Value = 1.83
$ws.Range("W10").Value = 10
$ws.Range("Y10").Value = 9.5
$ws.Range("Z10").Value = 8.5
$ws.Range("AA10").Value = 11
$ws.Range("AB10").Value = 23
$ws.Range("AC10").Value = 21
$ws.Range("AD10").Value = 13
$ws.Range("AE10").Value = 23
$ws.Range("AG10").Value = 600
$ws.Range("AH10").Value = 26
$ws.Range("AI10").Value = 41
$ws.Range("AJ10").Value = 23
$ws.Range("AK10").Value = 101
$ws.Range("AL10").Value = 51
$ws.Range("AN10").Value = 3.5
$ws.Range("AO10").Value = 5.5
$ws.Range("AP10").Value = 15
$ws.Range("AQ10").Value = 13
$ws.Range("AR10").Value = 34
$ws.Range("AS10").Value = 101
$ws.Range("AT10").Value = 4
$ws.Range("AX10").Value = 10
$ws.Range("AY10").Value = 41
$ws.Range("AZ10").Value = 41
$ws.Range("BA10").Value = 151
$ws.Range("BB10").Value = 151

# Row 11 updates
$ws.Range("A11").Value = "zgRmmLkg"
$ws.Range("C11").Value = "10:30"
$ws.Range("D11").Value = "SERBIA - SUPER LIGA"
$ws.Range("E11").Value = "Novi Pazar"
$ws.Range("F11").Value = "Vojvodina"
$ws.Range("G11").Value = 3.25
$ws.Range("H11").Value = 3.15
$ws.Range("I11").Value = 2.1
$ws.Range("J11").Value = 3.85
$ws.Range("K11").Value = 2.02
$ws.Range("L11").Value = 2.8
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 7.1
$ws.Range("O11").Value = 1.31
$ws.Range("P11").Value = 3.2
$ws.Range("Q11").Value = 1.91
$ws.Range("R11").Value = 1.8
$ws.Range("S11").Value = 1.44
$ws.Range("T11").Value = 2.6
$ws.Range("U11").Value = 1.72
$ws.Range("V11").Value = 2.02
$ws.Range("W11").Value = 10.25
$ws.Range("X11").Value = 18.5
$ws.Range("Y11").Value = 11.25
$ws.Range("Z11").Value = 45
$ws.Range("AA11").Value = 29
$ws.Range("AB11").Value = 35
$ws.Range("AC11").Value = 7.1
$ws.Range("AD11").Value = 6.2
$ws.Range("AE11").Value = 13
$ws.Range("AF11").Value = 60
$ws.Range("AG11").Value = 450
$ws.Range("AH11").Value = 7.5
$ws.Range("AI11").Value = 10.5
$ws.Range("AJ11").Value = 8.5
$ws.Range("AK11").Value = 21
$ws.Range("AL11").Value = 17
$ws.Range("AM11").Value = 27
$ws.Range("AN11").Value = 5.2
$ws.Range("AO11").Value = 18.5
$ws.Range("AP11").Value = 25
$ws.Range("AQ11").Value = 100
$ws.Range("AR11").Value = 120
$ws.Range("AS11").Value = 300
$ws.Range("AT11").Value = 2.6
$ws.Range("AU11").Value = 7
$ws.Range("AV11").Value = 65
$ws.Range("AW11").Value = ""
$ws.Range("AX11").Value = 4.05
$ws.Range("AY11").Value = 11.5
$ws.Range("AZ11").Value = 20
$ws.Range("BA11").Value = 45
$ws.Range("BB11").Value = 80
$ws.Range("BC11").Value = 250
$ws.Range("BD11").Value = ""

# Row 12 updates
$ws.Range("A12").Value = "tfYvksLt"
$ws.Range("C12").Value = "12:30"
$ws.Range("E12").Value = "Cukaricki"
$ws.Range("F12").Value = "Sp. Subotica"
$ws.Range("G12").Value = 1.6
$ws.Range("H12").Value = 3.55
$ws.Range("I12").Value = 5.3
$ws.Range("J12").Value = 2.15
$ws.Range("K12").Value = 2.15
$ws.Range("L12").Value = 5.6
$ws.Range("M12").Value = 1.08
$ws.Range("N12").Value = 6.7
$ws.Range("O12").Value = 1.37
$ws.Range("P12").Value = 2.87
$ws.Range("Q12").Value = 2.05
$ws.Range("R12").Value = 1.7
$ws.Range("S12").Value = 1.42
$ws.Range("T12").Value = 2.65
$ws.Range("U12").Value = 2.1
$ws.Range("V12").Value = 1.65
$ws.Range("W12").Value = 5.7
$ws.Range("X12").Value = 6.7
$ws.Range("Y12").Value = 8.5
$ws.Range("Z12").Value = 11.25
$ws.Range("AA12").Value = 14
$ws.Range("AC12").Value = 6.7
$ws.Range("AD12").Value = 7.2
$ws.Range("AE12").Value = 20
$ws.Range("AF12").Value = 120
$ws.Range("AG12").Value = 800
$ws.Range("AH12").Value = 12
$ws.Range("AI12").Value = 30
$ws.Range("AJ12").Value = 18.5
$ws.Range("AK12").Value = 110
$ws.Range("AL12").Value = 65
$ws.Range("AM12").Value = 75
$ws.Range("AN12").Value = 3.3
$ws.Range("AO12").Value = 7.6
$ws.Range("AP12").Value = 19
$ws.Range("AQ12").Value = 25
$ws.Range("AR12").Value = 65
$ws.Range("AT12").Value = 2.65
$ws.Range("AU12").Value = 8.25
$ws.Range("AV12").Value = 90
$ws.Range("AX12").Value = 6.9
$ws.Range("AY12").Value = 35
$ws.Range("AZ12").Value = 40
$ws.Range("BA12").Value = 250
$ws.Range("BB12").Value = 300
$ws.Range("BC12").Value = 500
